$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")

# --- Fix up the formatting of D43 so it matches the plain/default cell
# style used by the rest of the sheet (copy format only from a cell that
# already uses the plain "Normal" look, e.g. A2). ---
$ws.Range("A2").Copy()
$ws.Range("D43").PasteSpecial(-4122)

# --- Append four new rows (47-50) that mirror rows 43-46 (the PTIS
# Revenue/Commissioner approval chain) but represent the "...1" variants
# used for a second test data set. ---
$ws.Range("A47").Value = "PTISBillCollector1"
$ws.Range("B47").Value = "REVENUE"
$ws.Range("C47").Value = "Bill Collector"
$ws.Range("D47").Value = "PTBillCollectorOne~PTIS_REVBC_1"
$ws.Range("E47").Value = "Forward to bill collector"

$ws.Range("A48").Value = "PTISRevenueInspector1"
$ws.Range("B48").Value = "REVENUE"
$ws.Range("C48").Value = "UD Revenue Inspector"
$ws.Range("D48").Value = "PTRevenueInspectorOne~PTIS_REVSI_1"
$ws.Range("E48").Value = "Forward to revenue insoector"

$ws.Range("A49").Value = "PTISRevenueOfficer1"
$ws.Range("B49").Value = "REVENUE"
$ws.Range("C49").Value = "Revenue Officer"
$ws.Range("D49").Value = "PTRevenueOfficerOne~PTIS_REVOF_1"
$ws.Range("E49").Value = "Forward to revenue officer"

$ws.Range("A50").Value = "PTISCommissioner1"
$ws.Range("B50").Value = "ADMINISTRATION"
$ws.Range("C50").Value = "Commissioner"
$ws.Range("D50").Value = "PTCommissionerOne~ADM_COMM_20"

# --- Leave the cursor/selection the way the author left it. ---
[void]$ws.Range("B56").Select()
